$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (5th column), shifting existing
# columns E..L to F..M.
$ws.Columns("E:E").Insert()

# Set the header for the newly inserted column.
$ws.Range("E1").Value = "data_notificacio"

# Match the column width Excel applied for the new column (as seen in the
# target file: 15.140625).
$ws.Columns("E:E").ColumnWidth = 14.25

# Update the selected cell/range to match the post-edit workbook state.
$ws.Range("F1").Select()
